$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 0.29
$ws.Range("K2").Value = 30.3
$ws.Range("L2").Value = 0.3690621193666261
$ws.Range("U2").Value = 174.5
$ws.Range("V2").Value = 0.4665775401069518
$ws.Range("W2").Value = 0.0804140127388535
$ws.Range("X2").Value = 0.05374534097797576
$ws.Range("Y2").Value = 0.02666867176087774
$ws.Range("Z2").Value = 0.4383342231713828
$ws.Range("AB2").Value = 0.05374534097797576
$ws.Range("AC2").Value = -0.05374534097797576
$ws.Range("AG2").Value = -174.5
$ws.Range("AJ2").Value = -0.87468671679198
$ws.Range("AK2").Value = -0.7515073212747632

# Row 3 updates
$ws.Range("D3").Value = 0.29
$ws.Range("K3").Value = 30.3
$ws.Range("L3").Value = 0.3690621193666261
$ws.Range("U3").Value = 174.5
$ws.Range("V3").Value = 0.4665775401069518
$ws.Range("W3").Value = 0.0804140127388535
$ws.Range("X3").Value = 0.05374534097797576
$ws.Range("Y3").Value = 0.02666867176087774
$ws.Range("Z3").Value = 0.4383342231713828
$ws.Range("AB3").Value = 0.05374534097797576
$ws.Range("AC3").Value = -0.05374534097797576
$ws.Range("AG3").Value = -174.5
$ws.Range("AJ3").Value = -0.87468671679198
$ws.Range("AK3").Value = -0.7515073212747632
